$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.028.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "'2.174.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'250.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("D7").Value = "'67.05"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.53%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +3.66%  "
$ws.Range("D10").Value = "'36.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "'58.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'0.0929"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D14").Value = "'6.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.77%  "
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").Value = "'0.859"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "'14.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.40%  "
$ws.Range("D18").Value = "'2.175.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").Value = "'40.878.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0941"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Value = "'6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.99%  "
$ws.Range("D22").Value = "'71.49"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.17%  "
$ws.Range("D23").Value = "'230.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.50%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "'11.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +18.59%  "
$ws.Range("E26").Value = "  +5.92%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  -3.72%  "
$ws.Range("D30").Value = "'2.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").Value = "'169.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "'20.46"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("E33").Value = "  -1.68%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0721"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.27%  "
$ws.Range("D37").Value = "'4.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "'25.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.84%  "
$ws.Range("D39").Value = "'4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.26%  "
$ws.Range("D40").Value = "'0.0295"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.04%  "
$ws.Range("E41").Value = "  -2.72%  "
$ws.Range("D42").Value = "'12.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.91%  "
$ws.Range("D43").Value = "'5.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").Value = "'64.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'0.200"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.97%  "
$ws.Range("D46").Value = "'4.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.101"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'8.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.74%  "
$ws.Range("E49").Value = "  +0.36%  "
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("D51").Value = "'1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.56%  "
